$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the first two features as having no remaining hours (fully done).
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0

# Move the active selection to H7:I7, matching the user's next point of focus.
$ws.Range("H7:I7").Select()
